# Rename the three header labels that changed to match the samedis.care
# import template, and move the active selection to D2 (as in the diff).
#
# Column layout is unchanged; only these header cells' text changes:
#   C1: "Personalnummer"  -> "Mitarbeiternr."
#   D1: "Eintrittsdatum"  -> "Beitritt am"
#   E1: "Austrittsdatum"  -> "Austritt am"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mitarbeiter")

$ws.Range("C1").Value = "Mitarbeiternr."
$ws.Range("D1").Value = "Beitritt am"
$ws.Range("E1").Value = "Austritt am"

$ws.Activate()
$ws.Range("D2").Select()
